$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Add new row 63: a new "picd" (604) -> Diesel mapping entry, mirroring the
# existing rows' columns (picd, pids, Konto, Kontonavn, Prosess, Prosessnavn
# Posten, Konto tekst), so the "missing picd" check has a full mapping.
# Copy the formatting from the last existing data row (62) first, so the new
# row matches the table's look (bold/border on col A, plain border on the
# rest), then fill in the values.
$ws.Range("A62:G62").Copy()
$ws.Range("A63:G63").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A63").Value = 604
$ws.Range("B63").Value = "Diesel"
$ws.Range("C63").Value = 431100
$ws.Range("D63").Value = "Drivstoffkostnader, varekost"
$ws.Range("E63").Value = 1160
$ws.Range("F63").Value = " Diesel Innland"
$ws.Range("G63").Value = "Drivstoff"

# Grow Table1 (and its AutoFilter) to include the new row.
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:G63"))

# The _xlnm._FilterDatabase hidden name also needs to track the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='Ark1'!`$A`$1:`$G`$63"
    }
}

# Restore the view close to what the author ended up with: scrolled down to
# row 34 and the active cell on I53.
$ws.Range("I53").Select()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
